# edit.ps1
# Applies the "Updated symbol list" data refresh to the crypto tracking sheet.
# Coin/Link (columns B/C) are set as plain text; Price/Volume (columns D/E) are
# numeric- or percent-looking strings that must stay text, so they are entered
# with a leading apostrophe (Excel's text-force prefix) to avoid Excel silently
# converting them to numbers/percentages and losing exact formatting (e.g. trailing zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'245.70"
$ws.Range("E2").Value = "'1.19%"

# Row 3
$ws.Range("D3").Value = "'29.38"
$ws.Range("E3").Value = "'-1.39%"

# Row 4
$ws.Range("D4").Value = "'5.161"
$ws.Range("E4").Value = "'0.60%"

# Row 5
$ws.Range("D5").Value = "'0.05772"
$ws.Range("E5").Value = "'2.07%"

# Row 6
$ws.Range("D6").Value = "'6.599"
$ws.Range("E6").Value = "'1.64%"

# Row 7
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").Value = "'3.155"
$ws.Range("E7").Value = "'4.67%"

# Row 8
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = "'0.8585"
$ws.Range("E8").Value = "'3.76%"

# Row 9
$ws.Range("B9").Value = 'FTXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D9").Value = "'0.8619"
$ws.Range("E9").Value = "'-0.02%"

# Row 10
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = "'0.1363"
$ws.Range("E10").Value = "'2.52%"

# Row 11
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = "'0.07029"
$ws.Range("E11").Value = "'1.45%"

# Row 12
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D12").Value = "'0.03289"
$ws.Range("E12").Value = "'1.53%"

# Row 13
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = "'0.03018"
$ws.Range("E13").Value = "'5.60%"

# Row 14
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = "'0.09366"
$ws.Range("E14").Value = "'-0.24%"

# Row 15
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = "'0.001522"
$ws.Range("E15").Value = "'-0.07%"

# Row 16
$ws.Range("B16").Value = 'One'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D16").Value = "'0.0005978"
$ws.Range("E16").Value = "'-0.52%"

# Row 17
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").Value = "'0.006058"
$ws.Range("E17").Value = "'-1.78%"

# Row 18
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").Value = "'3.493"
$ws.Range("E18").Value = "'-0.79%"

# Row 19
$ws.Range("B19").Value = 'BTSEToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D19").Value = "'2.164"
$ws.Range("E19").Value = "'-2.42%"

# Row 20
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D20").Value = "'0.3201"
$ws.Range("E20").Value = "'1.66%"

# Row 21
$ws.Range("E21").Value = "'-1.00%"

# Row 22
$ws.Range("E22").Value = "'-8.57%"

# Row 23
$ws.Range("D23").Value = "'0.04135"
$ws.Range("E23").Value = "'-0.35%"

# Row 24
$ws.Range("D24").Value = "'0.1400"
$ws.Range("E24").Value = "'1.88%"

# Row 25
$ws.Range("E25").Value = "'1.34%"

# Row 26
$ws.Range("D26").Value = "'0.004136"
$ws.Range("E26").Value = "'-6.95%"

# Row 27
$ws.Range("E27").Value = "'2.59%"

# Row 28
$ws.Range("D28").Value = "'0.0001448"
$ws.Range("E28").Value = "'3.18%"

# Row 40
$ws.Range("E40").Value = "'0.65%"

# Row 41
$ws.Range("D41").Value = "'0.005889"
$ws.Range("E41").Value = "'72.64%"

# Row 42
$ws.Range("E42").Value = "'1.48%"

# Row 43
$ws.Range("E43").Value = "'-4.80%"

# Row 44
$ws.Range("D44").Value = "'0.008448"
$ws.Range("E44").Value = "'-12.48%"

# Row 45
$ws.Range("D45").Value = "'0.00005282"
$ws.Range("E45").Value = "'3.40%"

# Row 46
$ws.Range("E46").Value = "'-0.01%"

# Row 47
$ws.Range("D47").Value = "'0.05798"
$ws.Range("E47").Value = "'-42.58%"

# Row 48
$ws.Range("E48").Value = "'-4.44%"

# Row 49
$ws.Range("E49").Value = "'-0.01%"

# Row 50
$ws.Range("E50").Value = "'-0.01%"
